$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Pre-seed the shared-string table in the same order as the authored edit:
# AVR text (row 18) lands at index 24, HSERPRINT text (row 17) lands at index 25.
$ws.Cells.Item(18, 4).Value = "AVR compiles when it should not.  See https://sourceforge.net/p/gcbasic/discussion/596084/thread/e58866dc/#5e0f"
$ws.Cells.Item(17, 4).Value = "HSERPRINT not handling LONGs correctly.`nUpdated USART.H to handle LONGs correctly."

# Row 17 - Index 16, CLOSED, HSERPRINT LONGs fix (wrapped description style, like D2/D4/... cells)
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "CLOSED"

$ws.Cells.Item(17, 1).Style = "Normal"
$ws.Cells.Item(17, 2).Style = "Normal"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 1).HorizontalAlignment = -4131
$ws.Cells.Item(17, 1).VerticalAlignment = -4160
$ws.Cells.Item(17, 2).HorizontalAlignment = -4131
$ws.Cells.Item(17, 2).VerticalAlignment = -4160
$ws.Cells.Item(17, 4).HorizontalAlignment = -4131
$ws.Cells.Item(17, 4).VerticalAlignment = -4160
$ws.Cells.Item(17, 4).WrapText = $true
$ws.Rows.Item(17).RowHeight = 30

# Row 18 - Index 17, OPEN, AVR compiler bug (non-wrapped style, like A-col cells)
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "OPEN"

$ws.Cells.Item(18, 1).Style = "Normal"
$ws.Cells.Item(18, 2).Style = "Normal"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 1).HorizontalAlignment = -4131
$ws.Cells.Item(18, 1).VerticalAlignment = -4160
$ws.Cells.Item(18, 2).HorizontalAlignment = -4131
$ws.Cells.Item(18, 2).VerticalAlignment = -4160
$ws.Cells.Item(18, 4).HorizontalAlignment = -4131
$ws.Cells.Item(18, 4).VerticalAlignment = -4160
$ws.Cells.Item(18, 4).WrapText = $false

$ws.Range("D18").Select()

$ws.Application.ActiveWindow.ScrollRow = 13
